$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 18:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 537356
$ws.Range("C4").Value = 4477
$ws.Range("D4").Value = 31087
$ws.Range("E4").Value = 484834
$ws.Range("G4").Value = 858
$ws.Range("H4").Value = 21435

# Row 12 - Turquia
$ws.Range("B12").Value = 56956
$ws.Range("C12").Value = 4789
$ws.Range("D12").Value = 3446
$ws.Range("E12").Value = 52312
$ws.Range("F12").Value = 1665
$ws.Range("G12").Value = 97
$ws.Range("H12").Value = 1198

# Row 15 - Suiza
$ws.Range("B15").Value = 25407
$ws.Range("C15").Value = 300
$ws.Range("E15").Value = 12218

# Row 26 - Ecuador
$ws.Range("B26").Value = 7466
$ws.Range("C26").Value = 209
$ws.Range("D26").Value = 501
$ws.Range("E26").Value = 6632
$ws.Range("G26").Value = 18
$ws.Range("H26").Value = 333

# Row 43 - Serbia
$ws.Range("D43").Value = 400
$ws.Range("E43").Value = 3150
